$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains a new column L mirroring column K (the 2020 data column),
# representing a new "2021" year column appended to the table (rows 3-11).
# Column L3 is an empty, bottom-bordered cell like K3; L4 is the new year
# header (2021); L5:L10 repeat the same percentage figures as K5:K10; L11
# repeats the same total as K11. We copy K's formatting into L first, then
# set the values explicitly.
$ws.Range("K3:K11").Copy() | Out-Null
$ws.Range("L3:L11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("L4").Value = 2021
$ws.Range("L5").Value = 0.86
$ws.Range("L6").Value = 1.07
$ws.Range("L7").Value = 25.27
$ws.Range("L8").Value = 14
$ws.Range("L9").Value = 0.12
$ws.Range("L10").Value = 21.74
$ws.Range("L11").Value = 9.4600000000000009

# Move the active selection to N2, matching the post-edit view state.
$ws.Range("N2").Select() | Out-Null
